# Regenerate orders with updated distance/sizes.
# Distance codes: D51 -> D55, D64 -> D69, D80 -> D86
# Size code:      S30 -> S31
# These substrings appear throughout the Condition, Filename_Left,
# Filename_Right, Distance and Size columns (and their shared lookup
# values), so we do a global substring replace across the whole used
# range, same as Excel's Range.Replace (Ctrl+H) would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$xlPart = 2
$xlByRows = 1

$rng.Replace("D51", "D55", $xlPart, $xlByRows, $false, $false, $false, $false)
$rng.Replace("D64", "D69", $xlPart, $xlByRows, $false, $false, $false, $false)
$rng.Replace("D80", "D86", $xlPart, $xlByRows, $false, $false, $false, $false)
$rng.Replace("S30", "S31", $xlPart, $xlByRows, $false, $false, $false, $false)
